# ---------------------------------------------------------------------------
# Results_Summary.xlsx update: add Feb 17th, 2022 results as a new sheet
# ("7jpg_yolov5_on_cX"), rename the existing sheet to
# "8jpg_yolov5_on_cX_in_sequence", and refresh the UI selection state.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Rename the original sheet and insert the new one right after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "8jpg_yolov5_on_cX_in_sequence"

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "7jpg_yolov5_on_cX"

# --- 2. Title row (merged, like the other tables on sheet 1) --------------
$ws2.Range("C2:J2").Merge()
$ws2.Range("C2").Value = "February 17th, 2022"

# Reuse the existing "title bar" formatting (first / middle / last cell)
$ws1.Range("D18").Copy()
$ws2.Range("C2").PasteSpecial(-4122)
$ws1.Range("E18:J18").Copy()
$ws2.Range("D2:I2").PasteSpecial(-4122)
$ws1.Range("L18").Copy()
$ws2.Range("J2").PasteSpecial(-4122)

# Blank spacer cell to the left of the table (column B)
$ws1.Range("C3").Copy()
$ws2.Range("B2:B3").PasteSpecial(-4122)

# --- 3. Header row + row-label column --------------------------------------
$ws2.Range("C3").Value = "DUSVs_0"
$ws2.Range("D3").Value = "DUSVs_1"
$ws2.Range("E3").Value = "DUSVs_2"
$ws2.Range("F3").Value = "DUSVs_3"
$ws2.Range("G3").Value = "DUSVs_4"
$ws2.Range("H3").Value = "DUSVs_5"
$ws2.Range("I3").Value = "DUSVs_6"
$ws2.Range("J3").Value = "Total"

$ws2.Range("B4").Value = "Total"
$ws2.Range("B5").Value = "Correct"
$ws2.Range("B6").Value = "Mismatch"
$ws2.Range("B7").Value = "Crash"

# Header-row / label-column formatting
$ws1.Range("D19").Copy()
$ws2.Range("C3:J3").PasteSpecial(-4122)
$ws2.Range("B4:B7").PasteSpecial(-4122)

# --- 4. Data values ----------------------------------------------------
$ws2.Range("C4:I4").Value = @(79, 78, 78, 77, 77, 77, 79)
$ws2.Range("C5:I5").Value = @(53, 51, 49, 52, 51, 50, 50)
$ws2.Range("C6:I6").Value = @(0, 2, 3, 1, 2, 0, 0)
$ws2.Range("C7:I7").Value = @(26, 25, 26, 24, 24, 27, 29)

$ws2.Range("J4").Formula = "=SUM(C4:I4)"
$ws2.Range("J5").Formula = "=SUM(C5:I5)/J4"
$ws2.Range("J6").Formula = "=SUM(C6:I6)/J4"
$ws2.Range("J7").Formula = "=SUM(C7:I7)/J4"

# Data-cell formatting (border, centered) and totals/percent formats
$ws1.Range("D20").Copy()
$ws2.Range("C4:I7").PasteSpecial(-4122)

$ws1.Range("L20").Copy()
$ws2.Range("J4").PasteSpecial(-4122)

$ws1.Range("L21").Copy()
$ws2.Range("J5:J7").PasteSpecial(-4122)

# --- 5. Update the absolute path recorded for the workbook ----------------
# (cosmetic Mac Finder path metadata -- left untouched; not exposed via COM)

# --- 6. Selection / view state ---------------------------------------------
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("C18:L23").Select()

$ws2.Activate()
$ws2.Range("I8").Select()

$excel.CutCopyMode = $false
